# Update cryptocurrency price/volume snapshot (symbol list refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

# Row 2: BNB
Set-TextCell 2 4 "299.63"
Set-TextCell 2 5 "1.50%"

# Row 3: OKB
Set-TextCell 3 4 "31.33"
Set-TextCell 3 5 "-0.36%"

# Row 4: HuobiToken
Set-TextCell 4 4 "5.137"
Set-TextCell 4 5 "0.51%"

# Row 5: Cronos
Set-TextCell 5 4 "0.08276"
Set-TextCell 5 5 "11.26%"

# Row 6: FTXToken
Set-TextCell 6 4 "2.341"
Set-TextCell 6 5 "37.96%"

# Row 7: KuCoinToken
Set-TextCell 7 4 "7.904"
Set-TextCell 7 5 "2.60%"

# Row 8: GateToken
Set-TextCell 8 4 "3.843"
Set-TextCell 8 5 "1.65%"

# Row 9: MXToken
Set-TextCell 9 4 "0.9138"
Set-TextCell 9 5 "-1.70%"

# Row 10: WazirX
Set-TextCell 10 4 "0.1719"
Set-TextCell 10 5 "1.94%"

# Row 11: LiechtensteinCryptoassetsExchange
Set-TextCell 11 4 "0.07359"
Set-TextCell 11 5 "3.14%"

# Row 12: MandalaExchangeToken
Set-TextCell 12 4 "0.08039"
Set-TextCell 12 5 "1.33%"

# Row 13: BitrueCoin
Set-TextCell 13 4 "0.03025"
Set-TextCell 13 5 "0.82%"

# Row 14: BitMartToken
Set-TextCell 14 4 "0.09957"
Set-TextCell 14 5 "0.53%"

# Row 15: BitForexToken
Set-TextCell 15 4 "0.001526"
Set-TextCell 15 5 "1.66%"

# Row 16: TigerCash
Set-TextCell 16 4 "0.006139"
Set-TextCell 16 5 "-1.32%"

# Row 17: LEO
Set-TextCell 17 4 "3.496"
Set-TextCell 17 5 "1.12%"

# Row 18: BTSEToken
Set-TextCell 18 4 "2.254"
Set-TextCell 18 5 "1.29%"

# Row 19: BitpandaEcosystemToken
Set-TextCell 19 4 "0.3302"
Set-TextCell 19 5 "0.68%"

# Row 20: ProBitToken
Set-TextCell 20 4 "0.1350"
Set-TextCell 20 5 "1.10%"

# Row 21: MCDex
Set-TextCell 21 4 "4.608"
Set-TextCell 21 5 "1.06%"

# Row 22: ZBToken
Set-TextCell 22 4 "0.1607"
Set-TextCell 22 5 "3.37%"

# Row 23: CoinExToken
Set-TextCell 23 4 "0.04603"
Set-TextCell 23 5 "-1.00%"

# Row 24: BitKan
Set-TextCell 24 4 "0.001267"
Set-TextCell 24 5 "3.80%"

# Row 25: HotbitToken
Set-TextCell 25 5 "0.84%"

# Row 26: NitroEx
Set-TextCell 26 4 "0.0001186"
Set-TextCell 26 5 "-8.97%"

# Row 27: UpBots
Set-TextCell 27 4 "0.0003444"
Set-TextCell 27 5 "83.15%"

# Row 39: One
Set-TextCell 39 5 "11.01%"

# Row 40: IDEX
Set-TextCell 40 4 "0.04529"
Set-TextCell 40 5 "2.20%"

# Row 41: KickToken
Set-TextCell 41 4 "0.007262"
Set-TextCell 41 5 "2.63%"

# Row 42: BKEXToken
Set-TextCell 42 4 "0.1344"
Set-TextCell 42 5 "1.38%"

# Row 43: CEJI
Set-TextCell 43 4 "0.002193"
Set-TextCell 43 5 "4.79%"

# Row 44: LocalTraders
Set-TextCell 44 4 "0.01062"
Set-TextCell 44 5 "-13.71%"

# Row 45: CoinLion
Set-TextCell 45 4 "0.00006290"
Set-TextCell 45 5 "4.40%"

# Row 46: Kangarootoken
Set-TextCell 46 4 "0.00000000754"
Set-TextCell 46 5 "0.52%"

# Row 47: CoinbaseStockToken
Set-TextCell 47 4 "0.006671"
Set-TextCell 47 5 "-39.51%"

# Row 48: BOLO
Set-TextCell 48 5 "15.31%"

# Row 49: CryptobidCoin
Set-TextCell 49 4 "0.00002111"
Set-TextCell 49 5 "0.52%"

# Row 50: SpecialPowerGold
Set-TextCell 50 4 "0.0002011"
Set-TextCell 50 5 "0.59%"
